$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.924.56'
$ws.Range('E2').Value = '  -0.26%  '
$ws.Range('D3').Value = '1.622.66'
$ws.Range('E3').Value = '  -1.15%  '
$ws.Range('E4').Value = '  +0.16%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '213.63'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('E8').Value = '  -2.27%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0616'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -3.45%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.33'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -6.12%  '
$ws.Range('E11').Value = '  -1.13%  '
$ws.Range('D12').Value = '1.848.99'
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('B13').Value = 'Polkadot'
$ws.Range('C13').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '4.18'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  -2.13%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').Value = '1.615.39'
$ws.Range('E14').Value = '  -1.83%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.525'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -3.51%  '
$ws.Range('D16').Value = '25.909.71'
$ws.Range('E16').Value = '  -0.42%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '61.15'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -3.52%  '
$ws.Range('D18').Value = '0.0₃0733'
$ws.Range('E18').Value = '  -3.85%  '
$ws.Range('E19').Value = '  +0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.90'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -1.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.23'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -3.00%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.59'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -3.35%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '6.07'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -2.13%  '
$ws.Range('E24').Value = '  +0.26%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.57'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.39%  '
$ws.Range('E26').Value = '  +0.20%  '
$ws.Range('E27').Value = '  -3.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.71'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -2.35%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.13'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  -2.50%  '
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0482'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.70%  '
$ws.Range('E32').Value = '  -4.97%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.10'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -5.71%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.50'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -2.62%  '
$ws.Range('E35').Value = '  -2.36%  '
$ws.Range('D36').Value = '1.118.25'
$ws.Range('E36').Value = '  -1.02%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.847'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -6.35%  '
$ws.Range('E38').Value = '  -1.54%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.517'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -4.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0153'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -2.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '97.96'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.08%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.765'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -4.05%  '
$ws.Range('B43').Value = 'RocketPoolETH'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D43').Value = '1.758.05'
$ws.Range('E43').Value = '  -1.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.14'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -5.85%  '
$ws.Range('E45').Value = '  -1.69%  '
$ws.Range('E46').Value = '  +1.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '54.29'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.03%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.47'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.45%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.412'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.40%  '
$ws.Range('E50').Value = '  +0.38%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '7.48'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -3.39%  '
